$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the old selected-tab marker on Sheet1 (selection/view only; handled automatically
# once a later sheet becomes active, but set explicitly for clarity).

# --- Insert the two new sheets, in order, right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Nano Degree"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Front-End Master"

# --- Column widths (character units; engine quantizes to 1/6-char steps) ---
$ws2.Columns.Item(2).ColumnWidth = 87.92
$ws2.Columns.Item(3).ColumnWidth = 79.42203125
$ws3.Columns.Item(2).ColumnWidth = 87.7540625
$ws3.Columns.Item(3).ColumnWidth = 79.259375

# --- Style A: font Arial, vertical-center + wrap (no horizontal align) ---
# --- Style B: font Arial only (used for blank spacer cells) ---
# --- Style C: font Arial colored #1F2932, left+vertical-center + wrap ---

# --- Seed style A on Front-End Master!B2, style B on Front-End Master!B7,
#     then style C on Nano Degree!B2 -- ORDER MATTERS for a clean style table. ---
$ws3.Range("B2").Font.Name = "Arial"
$ws3.Range("B2").VerticalAlignment = -4108
$ws3.Range("B2").WrapText = $true
$ws3.Range("B7").Font.Name = "Arial"
$ws2.Range("B2").HorizontalAlignment = -4131
$ws2.Range("B2").VerticalAlignment = -4108
$ws2.Range("B2").WrapText = $true
$ws2.Range("B2").Font.Name = "Arial"
$ws2.Range("B2").Font.Color = 3287327

# --- Nano Degree sheet values ---
$ws2.Range("B2").Value = 'Front End Web Developer  '
$ws2.Range("C2").Value = 'https://classroom.udacity.com/nanodegrees/nd0011/dashboard/overview'
$ws2.Range("B6").Value = '26 module'

# --- Front-End Master sheet values + styles ---
$ws3.Range("B2").Value = 'Complete Intro to Web Development, v2'
$ws3.Range("C2").Value = 'https://frontendmasters.com/courses/web-development-v2/'
$ws3.Range("B4").Value = 'CSS In-Depth, v2'
$ws3.Range("B4").Font.Name = "Arial"
$ws3.Range("B4").VerticalAlignment = -4108
$ws3.Range("B4").WrapText = $true
$ws3.Range("C4").Value = 'https://frontendmasters.com/courses/css-in-depth-v2/'
$ws3.Range("B6").Value = 'CSS Grids and Flexbox for Responsive Web Design'
$ws3.Range("B6").Font.Name = "Arial"
$ws3.Range("B6").VerticalAlignment = -4108
$ws3.Range("B6").WrapText = $true
$ws3.Range("C6").Value = 'https://frontendmasters.com/courses/css-grids-flexbox/'
$ws3.Range("B8").Value = 'Advanced CSS Layouts'
$ws3.Range("B8").Font.Name = "Arial"
$ws3.Range("B8").VerticalAlignment = -4108
$ws3.Range("B8").WrapText = $true
$ws3.Range("C8").Value = 'https://frontendmasters.com/courses/advanced-css-layouts/'
$ws3.Range("B10").Value = 'Webpack 4 Fundamentals'
$ws3.Range("B10").Font.Name = "Arial"
$ws3.Range("B10").VerticalAlignment = -4108
$ws3.Range("B10").WrapText = $true
$ws3.Range("C10").Value = 'https://frontendmasters.com/courses/webpack-fundamentals/'
$ws3.Range("B11").Font.Name = "Arial"
$ws3.Range("B12").Value = 'Web Performance with Webpack'
$ws3.Range("B12").Font.Name = "Arial"
$ws3.Range("B12").VerticalAlignment = -4108
$ws3.Range("B12").WrapText = $true
$ws3.Range("C12").Value = 'https://frontendmasters.com/courses/performance-webpack/'
$ws3.Range("B13").Font.Name = "Arial"
$ws3.Range("B14").Value = 'Webpack Plugins System'
$ws3.Range("B14").Font.Name = "Arial"
$ws3.Range("B14").VerticalAlignment = -4108
$ws3.Range("B14").WrapText = $true
$ws3.Range("C14").Value = 'https://frontendmasters.com/courses/webpack-plugins/'
$ws3.Range("B16").Value = 'Sass Fundamentals'
$ws3.Range("B16").Font.Name = "Arial"
$ws3.Range("B16").VerticalAlignment = -4108
$ws3.Range("B16").WrapText = $true
$ws3.Range("C16").Value = 'https://frontendmasters.com/courses/sass/'
$ws3.Range("B17").Font.Name = "Arial"
$ws3.Range("B18").Value = 'Secure Authentication for Web Apps & APIs Using JWTs'
$ws3.Range("B18").Font.Name = "Arial"
$ws3.Range("B18").VerticalAlignment = -4108
$ws3.Range("B18").WrapText = $true
$ws3.Range("C18").Value = 'https://frontendmasters.com/courses/secure-auth-jwt/'
$ws3.Range("B19").Font.Name = "Arial"
$ws3.Range("B20").Value = 'Digging Into Node.js'
$ws3.Range("B20").Font.Name = "Arial"
$ws3.Range("B20").VerticalAlignment = -4108
$ws3.Range("B20").WrapText = $true
$ws3.Range("C20").Value = 'https://frontendmasters.com/courses/digging-into-node/'
$ws3.Range("B22").Value = 'Introduction to Dev Tools, v3'
$ws3.Range("B22").Font.Name = "Arial"
$ws3.Range("B22").VerticalAlignment = -4108
$ws3.Range("B22").WrapText = $true
$ws3.Range("C22").Value = 'https://frontendmasters.com/workshops/dev-tools-v3/'
$ws3.Range("B23").Font.Name = "Arial"
$ws3.Range("B24").Value = 'JavaScript Performance'
$ws3.Range("B24").Font.Name = "Arial"
$ws3.Range("B24").VerticalAlignment = -4108
$ws3.Range("B24").WrapText = $true
$ws3.Range("C24").Value = 'https://frontendmasters.com/courses/web-performance/'
$ws3.Range("B25").Font.Name = "Arial"
$ws3.Range("B26").Value = 'TypeScript 3 Fundamentals, v2'
$ws3.Range("B26").Font.Name = "Arial"
$ws3.Range("B26").VerticalAlignment = -4108
$ws3.Range("B26").WrapText = $true
$ws3.Range("C26").Value = 'https://frontendmasters.com/courses/typescript-v2/'
$ws3.Range("B28").Value = 'JavaScript: The Hard Parts, v2'
$ws3.Range("B28").Font.Name = "Arial"
$ws3.Range("B28").VerticalAlignment = -4108
$ws3.Range("B28").WrapText = $true
$ws3.Range("C28").Value = 'https://frontendmasters.com/courses/javascript-hard-parts-v2/'
$ws3.Range("B30").Value = 'Responsive Web Typography v2'
$ws3.Range("B30").Font.Name = "Arial"
$ws3.Range("B30").VerticalAlignment = -4108
$ws3.Range("B30").WrapText = $true
$ws3.Range("C30").Value = 'https://frontendmasters.com/courses/responsive-typography-v2/'
$ws3.Range("B31").Font.Name = "Arial"
$ws3.Range("B32").Value = 'HTML Email Development, v2'
$ws3.Range("B32").Font.Name = "Arial"
$ws3.Range("B32").VerticalAlignment = -4108
$ws3.Range("B32").WrapText = $true
$ws3.Range("C32").Value = 'https://frontendmasters.com/courses/html-email-v2/'
$ws3.Range("B34").Value = 'Introduction to Next.js'
$ws3.Range("B34").Font.Name = "Arial"
$ws3.Range("B34").VerticalAlignment = -4108
$ws3.Range("B34").WrapText = $true
$ws3.Range("C34").Value = 'https://frontendmasters.com/courses/next-js/'
$ws3.Range("B35").Font.Name = "Arial"
$ws3.Range("B36").Value = 'Production-Grade Next.js'
$ws3.Range("B36").Font.Name = "Arial"
$ws3.Range("B36").VerticalAlignment = -4108
$ws3.Range("B36").WrapText = $true
$ws3.Range("C36").Value = 'https://frontendmasters.com/courses/production-next/'
$ws3.Range("B38").Value = 'Modern Search Engine Optimization (SEO)'
$ws3.Range("B38").Font.Name = "Arial"
$ws3.Range("B38").VerticalAlignment = -4108
$ws3.Range("B38").WrapText = $true
$ws3.Range("C38").Value = 'https://frontendmasters.com/courses/modern-seo/'

# --- Sheet view / selection state ---
$ws1.Range("C26").Select()
$ws2.Range("B6").Select()
$ws3.Activate()
$ws3.Range("C40").Select()

# --- Page setup for Front-End Master (portrait, A4-ish paper 9) ---
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

Write-Output "edit complete"
